$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titles = @(
    "Vice President Vance meets Zelenskyy, Pearson airport gold heist, GST holiday and more",
    "Bockstael Construction Commits `$100,000 to the HSC Foundation’s Operation Excellence Corporate Challenge",
    "ROSEN, TOP RANKED GLOBAL COUNSEL, Encourages Grocery Outlet Holding Corp. Investors to Secure Counsel Before Important Deadline in Securities Class Action - GO",
    "Sources: Paramount execs weigh risks of settling Trump's lawsuit against CBS News, concerned about potential shareholder litigation or criminal bribery charges (Jessica Toonkel/Wall Street Journal)",
    "BREAKING: SEC Says Lawsuit Against Coinbase May End - Hot Moments",
    "ROSEN, A LEADING LAW FIRM, Encourages FTAI Aviation Ltd. Investors to Secure Counsel Before Important Deadline in Securities Class Action - FTAI",
    "ROSEN, A TOP RANKED LAW FIRM, Encourages Innovative Industrial Properties, Inc. Investors to Secure Counsel Before Important Deadline in Securities Class Action — IIPR",
    "NVO STOCK NEWS: NVO Shareholders with Large Losses Should Contact Robbins LLP for Information About the Class Action Lawsuit Against Novo Nordisk A/S",
    "Natasha accused of allegedly forcing 2Face to post divorce video",
    "WANTED FOR THEFT IN HUNTSVILLE"
)

$startRow = 145
for ($i = 0; $i -lt $titles.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $titles[$i]
}
